# Update "想去人数" (want-to-go count) figures across sheets, per
# the latest generated-output refresh (gh-pages @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3212
$ws1.Range("F3").Value = 734
$ws1.Range("F5").Value = 6880
$ws1.Range("F6").Value = 2053
$ws1.Range("F7").Value = 25
$ws1.Range("F12").Value = 23
$ws1.Range("F13").Value = 149
$ws1.Range("F14").Value = 188

# --- 演出 (Performances) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 17

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3212
$ws4.Range("F3").Value = 17
$ws4.Range("F4").Value = 734
$ws4.Range("F6").Value = 6880
$ws4.Range("F7").Value = 2053
$ws4.Range("F8").Value = 25
$ws4.Range("F13").Value = 23
$ws4.Range("F14").Value = 149
$ws4.Range("F15").Value = 188
